$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the Sayon Das / SDET row with Konul / X / Tester,
# keeping Elkhan / N / Automation Tester as the row below it.
$ws.Range("A2").Value = "Konul"
$ws.Range("B2").Value = "X"
$ws.Range("C2").Value = "Tester"

$ws.Range("A3").Value = "Elkhan"
$ws.Range("B3").Value = "N"
$ws.Range("C3").Value = "Automation Tester"

# Update the active selection on the sheet.
$ws.Range("C6").Select()
